$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: GROUP BY + AGGREGATE -> AGGREGATION
$ws.Range("A10").Value = "AGGREGATION"

# Row 22: LARGE COMPLEX -> LARGE
$ws.Range("A22").Value = "LARGE"

# Row 26: COUNT(*) ROWS -> COUNT(*), B26 rename
$ws.Range("A26").Value = "COUNT(*)"
$ws.Range("B26").Value = "manual_test_agg_all_3"

# Row 27: B27 rename
$ws.Range("B27").Value = "high_level_test_agg_all_3"

# Row 28: B28 rename
$ws.Range("B28").Value = "manual_test_agg_all_5"

# Row 29: B29 rename
$ws.Range("B29").Value = "high_level_test_agg_all_5"
